$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting the existing
# N/O ("Late" header + values) and O/P ("Outstanding" header + values)
# one column to the right.
$ws.Columns("N").Insert()

# Excel assigns a plain (non best-fit) width to a freshly inserted column;
# match the authored width of 10 characters.
$ws.Columns("N").ColumnWidth = 9.166666666666666

# Restore the active selection recorded after the edit.
[void]$ws.Range("R9").Select()
